# edit.ps1 — apply the "User Stories" formatting/content changes:
#  1) Title block ("FIT 3077" / "Semester 1" / "User Stories") switches
#     from hard-coded Arial/Times New Roman fonts to the theme fonts
#     (majorHAnsi), and the "User"/"ser" run split is merged into one
#     "User" run.
#  2) The placeholder group-name heading text becomes the real group name.
#  3) The broken-up, spell-check-annotated author-name runs are merged
#     back into a single run.
#
# Because this host's simplified Font object doesn't expose the OOXML
# "theme font" attributes (asciiTheme/hAnsiTheme/cstheme), we set the
# precise run/paragraph XML directly via Range.InsertXML so the
# <w:rFonts .../> element ends up exactly as intended.

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $rng = $paragraph.Range
    # Exclude the trailing paragraph mark from the replaced range so we
    # don't disturb the following paragraph; InsertXML replaces the
    # range's content (including supplying its own paragraph mark from
    # the <w:p> element we pass in).
    $target = $d.Range($rng.Start, $rng.End)
    $target.InsertXML($innerXml)
}

function Find-ParagraphByText([string]$needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like ($needle + "*")) {
            return $p
        }
    }
    return $null
}

$wNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14NS = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------
# 1) "FIT 3077" paragraph: ascii/hAnsi/cs -> theme majorHAnsi fonts.
# ---------------------------------------------------------------------
$p1 = Find-ParagraphByText "FIT 3077"
if ($null -eq $p1) { throw "Could not locate the 'FIT 3077' paragraph" }
$xml1 = '<w:p ' + $wNS + ' ' + $w14NS + ' w14:paraId="1BEAE7DD" w14:textId="41245913" w:rsidR="004C0FA8" w:rsidRPr="004C0FA8" w:rsidRDefault="004C0FA8" w:rsidP="004C0FA8">' + `
    '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:jc w:val="center"/>' + `
    '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="004C0FA8"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>FIT 3077</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml $p1 $xml1

# ---------------------------------------------------------------------
# 2) "Semester 1" paragraph: ascii/hAnsi/cs -> theme majorHAnsi fonts
#    (eastAsia stays Times New Roman).
# ---------------------------------------------------------------------
$p2 = Find-ParagraphByText "Semester 1"
if ($null -eq $p2) { throw "Could not locate the 'Semester 1' paragraph" }
$xml2 = '<w:p ' + $wNS + ' ' + $w14NS + ' w14:paraId="1DCDC3A8" w14:textId="77777777" w:rsidR="004C0FA8" w:rsidRPr="004C0FA8" w:rsidRDefault="004C0FA8" w:rsidP="004C0FA8">' + `
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="004C0FA8"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Semester 1</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml $p2 $xml2

# ---------------------------------------------------------------------
# 3) "User Stories" title paragraph: fonts -> theme majorHAnsi, and the
#    "U" + "ser" runs merge into a single "User" run.
# ---------------------------------------------------------------------
$p3 = Find-ParagraphByText "User Stories"
if ($null -eq $p3) { throw "Could not locate the 'User Stories' title paragraph" }
$xml3 = '<w:p ' + $wNS + ' ' + $w14NS + ' w14:paraId="6684C6EB" w14:textId="19ECC5DB" w:rsidR="004C0FA8" w:rsidRPr="004C0FA8" w:rsidRDefault="004C0FA8" w:rsidP="004C0FA8">' + `
    '<w:pPr><w:spacing w:after="60"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>User</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Stories</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml $p3 $xml3

# ---------------------------------------------------------------------
# 4) Group-name heading: replace the placeholder text with the group
#    name (formatting unchanged).
# ---------------------------------------------------------------------
$p4 = Find-ParagraphByText ([char]0x201C + "Group Name" + [char]0x201D + ":")
if ($null -eq $p4) { throw "Could not locate the group-name placeholder paragraph" }
$xml4 = '<w:p ' + $wNS + ' ' + $w14NS + ' w14:paraId="5FEE24B7" w14:textId="77777777" w:rsidR="004C0FA8" w:rsidRPr="004C0FA8" w:rsidRDefault="004C0FA8" w:rsidP="004C0FA8">' + `
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="004C0FA8"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Torino Development United</w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml $p4 $xml4

# ---------------------------------------------------------------------
# 5) Author names paragraph: merge the spell-check-split runs
#    (with <w:proofErr/> markers) back into one run.
# ---------------------------------------------------------------------
$p5 = Find-ParagraphByText "Soo Guan Yin"
if ($null -eq $p5) { throw "Could not locate the author-names paragraph" }
$xml5 = '<w:p ' + $wNS + ' ' + $w14NS + ' w14:paraId="4E035ED1" w14:textId="77777777" w:rsidR="004C0FA8" w:rsidRPr="004C0FA8" w:rsidRDefault="004C0FA8" w:rsidP="004C0FA8">' + `
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="004C0FA8"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Soo Guan Yin, Chua Jun Jie, Justin Chuah, Lim Fluoryynx </w:t></w:r>' + `
    '</w:p>'
Set-ParagraphXml $p5 $xml5

Write-Output "done"
